# Added the TLC59711 everywhere!! -- add the missing part row to the micro BOM (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM line: U102 / TLC59711 (LED driver) that had been omitted.
# Columns: A Quantity | B Refdes | C Value | D Footprint/Package |
#          E Specification | F Supplier | G Supplier PN | H Manufacturer |
#          I Manufacturer PN | J Price in $ (10/100/500)
$ws.Range("A25").Value = 1
$ws.Range("B25").Value = "U102"
$ws.Range("C25").Value = "TLC59711"
$ws.Range("D25").Value = "HTSSOP-20"
$ws.Range("F25").Value = "Mouser"
$ws.Range("G25").Value = "595-TLC59711PWP"
$ws.Range("H25").Value = "TI"
$ws.Range("I25").Value = "TLC59711PWP"
$ws.Range("J25").Value = "4.12/3.37/2.47"

# Leave the selection where the author left it when they finished editing.
$ws.Range("F26").Select()
